$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44363
$ws.Range("J2").Value = 250
$ws.Range("K2").Value = 2500
$ws.Range("L2").Value = 2800
$ws.Range("M2").Value = 2650
$ws.Range("N2").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("P2").Value = 1325
$ws.Range("Q2").Value = 2

$ws.Range("D3").Value = 44917
$ws.Range("J3").Value = 300
$ws.Range("K3").Value = 2700
$ws.Range("L3").Value = 3000
$ws.Range("M3").Value = 2850
$ws.Range("N3").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("P3").Value = 1425
$ws.Range("Q3").Value = 2

$ws.Range("D4").Value = 44438
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 950
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = 975
$ws.Range("N4").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("P4").Value = 488
$ws.Range("Q4").Value = 2

$ws.Range("D5").Value = 44365
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 1800
$ws.Range("L5").Value = 2000
$ws.Range("M5").Value = 1900
$ws.Range("N5").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("P5").Value = 950
$ws.Range("Q5").Value = 2

$ws.Range("D6").Value = 44468
$ws.Range("J6").Value = 300
$ws.Range("K6").Value = 900
$ws.Range("L6").Value = 1000
$ws.Range("M6").Value = 950
$ws.Range("N6").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("P6").Value = 475
$ws.Range("Q6").Value = 2

$ws.Range("D7").Value = 44726
$ws.Range("J7").Value = 250
$ws.Range("K7").Value = 2500
$ws.Range("L7").Value = 2800
$ws.Range("M7").Value = 2650
$ws.Range("N7").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("P7").Value = 1325
$ws.Range("Q7").Value = 2

$ws.Range("D8").Value = 44302
$ws.Range("J8").Value = 300
$ws.Range("K8").Value = 900
$ws.Range("L8").Value = 1000
$ws.Range("M8").Value = 950
$ws.Range("N8").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("P8").Value = 475
$ws.Range("Q8").Value = 2

$ws.Range("D9").Value = 44971
$ws.Range("J9").Value = 350
$ws.Range("K9").Value = 2500
$ws.Range("L9").Value = 2800
$ws.Range("M9").Value = 2671
$ws.Range("N9").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("P9").Value = 1336
$ws.Range("Q9").Value = 2

$ws.Range("D10").Value = 44172
$ws.Range("J10").Value = 200
$ws.Range("K10").Value = 1300
$ws.Range("L10").Value = 1500
$ws.Range("M10").Value = 1400
$ws.Range("N10").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("P10").Value = 700
$ws.Range("Q10").Value = 2

$ws.Range("D11").Value = 44253
$ws.Range("J11").Value = 250
$ws.Range("K11").Value = 1800
$ws.Range("L11").Value = 2000
$ws.Range("M11").Value = 1900
$ws.Range("N11").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("P11").Value = 950
$ws.Range("Q11").Value = 2

$ws.Range("D12").Value = 44427
$ws.Range("J12").Value = 250
$ws.Range("K12").Value = 1300
$ws.Range("L12").Value = 1500
$ws.Range("M12").Value = 1400
$ws.Range("N12").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("P12").Value = 700
$ws.Range("Q12").Value = 2

$ws.Range("D13").Value = 44435
$ws.Range("J13").Value = 300
$ws.Range("K13").Value = 900
$ws.Range("L13").Value = 1000
$ws.Range("M13").Value = 950
$ws.Range("N13").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("P13").Value = 475
$ws.Range("Q13").Value = 2

$ws.Range("D14").Value = 44229
$ws.Range("J14").Value = 250
$ws.Range("K14").Value = 1800
$ws.Range("L14").Value = 2000
$ws.Range("M14").Value = 1900
$ws.Range("N14").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("P14").Value = 950
$ws.Range("Q14").Value = 2

$ws.Range("D15").Value = 44572
$ws.Range("J15").Value = 300
$ws.Range("K15").Value = 1400
$ws.Range("L15").Value = 1500
$ws.Range("M15").Value = 1450
$ws.Range("N15").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("P15").Value = 725
$ws.Range("Q15").Value = 2

$ws.Range("D16").Value = 44390
$ws.Range("J16").Value = 250
$ws.Range("K16").Value = 2400
$ws.Range("L16").Value = 2500
$ws.Range("M16").Value = 2450
$ws.Range("N16").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("P16").Value = 1225
$ws.Range("Q16").Value = 2

$ws.Range("D17").Value = 44181
$ws.Range("J17").Value = 200
$ws.Range("K17").Value = 1000
$ws.Range("L17").Value = 1200
$ws.Range("M17").Value = 1100
$ws.Range("N17").Value = '$/atado'
$ws.Range("P17").Value = 1100
$ws.Range("Q17").Value = 1

$ws.Range("D18").Value = 44544
$ws.Range("J18").Value = 250
$ws.Range("K18").Value = 900
$ws.Range("L18").Value = 1000
$ws.Range("M18").Value = 950
$ws.Range("N18").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("P18").Value = 475
$ws.Range("Q18").Value = 2

$ws.Range("D19").Value = 44202
$ws.Range("J19").Value = 250
$ws.Range("K19").Value = 1800
$ws.Range("L19").Value = 2000
$ws.Range("M19").Value = 1900
$ws.Range("N19").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("P19").Value = 950
$ws.Range("Q19").Value = 2

$ws.Range("D20").Value = 44243
$ws.Range("J20").Value = 250
$ws.Range("K20").Value = 1200
$ws.Range("L20").Value = 1300
$ws.Range("M20").Value = 1250
$ws.Range("N20").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("P20").Value = 625
$ws.Range("Q20").Value = 2

$ws.Range("D21").Value = 44795
$ws.Range("J21").Value = 250
$ws.Range("K21").Value = 1800
$ws.Range("L21").Value = 2000
$ws.Range("M21").Value = 1900
$ws.Range("N21").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("P21").Value = 950
$ws.Range("Q21").Value = 2

$ws.Range("D22").Value = 44266
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 1700
$ws.Range("L22").Value = 1800
$ws.Range("M22").Value = 1750
$ws.Range("N22").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("P22").Value = 875
$ws.Range("Q22").Value = 2

$ws.Range("D23").Value = 44936
$ws.Range("J23").Value = 350
$ws.Range("K23").Value = 3000
$ws.Range("L23").Value = 3500
$ws.Range("M23").Value = 3357
$ws.Range("N23").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("P23").Value = 1678
$ws.Range("Q23").Value = 2

$ws.Range("D24").Value = 44972
$ws.Range("J24").Value = 350
$ws.Range("K24").Value = 800
$ws.Range("L24").Value = 1000
$ws.Range("M24").Value = 943
$ws.Range("N24").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("P24").Value = 472
$ws.Range("Q24").Value = 2

$ws.Range("D25").Value = 44817
$ws.Range("J25").Value = 300
$ws.Range("K25").Value = 900
$ws.Range("L25").Value = 1000
$ws.Range("M25").Value = 950
$ws.Range("N25").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("P25").Value = 475
$ws.Range("Q25").Value = 2

$ws.Range("D26").Value = 44447
$ws.Range("J26").Value = 300
$ws.Range("K26").Value = 900
$ws.Range("L26").Value = 1000
$ws.Range("M26").Value = 950
$ws.Range("N26").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("P26").Value = 475
$ws.Range("Q26").Value = 2

$ws.Range("D27").Value = 44601
$ws.Range("J27").Value = 270
$ws.Range("K27").Value = 2200
$ws.Range("L27").Value = 2500
$ws.Range("M27").Value = 2350
$ws.Range("N27").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("P27").Value = 1175
$ws.Range("Q27").Value = 2

$ws.Range("D28").Value = 44540
$ws.Range("J28").Value = 300
$ws.Range("K28").Value = 900
$ws.Range("L28").Value = 1000
$ws.Range("M28").Value = 950
$ws.Range("N28").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("P28").Value = 475
$ws.Range("Q28").Value = 2

$ws.Range("D29").Value = 44257
$ws.Range("J29").Value = 500
$ws.Range("K29").Value = 1400
$ws.Range("L29").Value = 1500
$ws.Range("M29").Value = 1450
$ws.Range("N29").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("P29").Value = 725
$ws.Range("Q29").Value = 2

$ws.Range("D30").Value = 44525
$ws.Range("J30").Value = 300
$ws.Range("K30").Value = 1400
$ws.Range("L30").Value = 1500
$ws.Range("M30").Value = 1450
$ws.Range("N30").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("P30").Value = 725
$ws.Range("Q30").Value = 2

$ws.Range("D31").Value = 44789
$ws.Range("J31").Value = 300
$ws.Range("K31").Value = 1400
$ws.Range("L31").Value = 1500
$ws.Range("M31").Value = 1450
$ws.Range("N31").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("P31").Value = 725
$ws.Range("Q31").Value = 2

$ws.Range("D32").Value = 44403
$ws.Range("J32").Value = 250
$ws.Range("K32").Value = 1800
$ws.Range("L32").Value = 2000
$ws.Range("M32").Value = 1900
$ws.Range("N32").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("P32").Value = 950
$ws.Range("Q32").Value = 2

$ws.Range("D33").Value = 44616
$ws.Range("J33").Value = 270
$ws.Range("K33").Value = 1300
$ws.Range("L33").Value = 1500
$ws.Range("M33").Value = 1400
$ws.Range("N33").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("P33").Value = 700
$ws.Range("Q33").Value = 2

$ws.Range("D34").Value = 44385
$ws.Range("J34").Value = 300
$ws.Range("K34").Value = 2400
$ws.Range("L34").Value = 2500
$ws.Range("M34").Value = 2450
$ws.Range("N34").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("P34").Value = 1225
$ws.Range("Q34").Value = 2

$ws.Range("D35").Value = 44291
$ws.Range("J35").Value = 250
$ws.Range("K35").Value = 1800
$ws.Range("L35").Value = 2000
$ws.Range("M35").Value = 1900
$ws.Range("N35").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("P35").Value = 950
$ws.Range("Q35").Value = 2

$ws.Range("D36").Value = 44392
$ws.Range("J36").Value = 250
$ws.Range("K36").Value = 1800
$ws.Range("L36").Value = 2000
$ws.Range("M36").Value = 1900
$ws.Range("N36").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("P36").Value = 950
$ws.Range("Q36").Value = 2

$ws.Range("D37").Value = 44161
$ws.Range("J37").Value = 270
$ws.Range("K37").Value = 900
$ws.Range("L37").Value = 1000
$ws.Range("M37").Value = 950
$ws.Range("N37").Value = '$/atado 1,5 a 2 kilos'
$ws.Range("P37").Value = 475
$ws.Range("Q37").Value = 2
